# HTML GRAPH Part4 - Inprogress
# Adds the latest regression run rows to the AMSIN and BETA sprint-history
# sheets, and refreshes the "Run Time" timestamp cell format/value.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd h:mm:ss"

function Set-TextValue($cell, [string]$text) {
    # Force the cell to store a literal text value instead of letting Excel
    # auto-convert date-like / numeric-looking strings, then drop back to the
    # default (unstyled) cell format so no stray formatting is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# AMSIN sheet: correct the existing B2 run-time value/format and append the
# 2021-06-07 regression run as row 3.
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Create a throwaway cell carrying the "yyyy-mm-dd h:mm:ss" format first, so
# that when it's applied to B2/B3 below they pick up a *new* shared style
# entry instead of reusing B2's current one (matching how the workbook's
# style table grows when the format is (re)applied from scratch).
$scratch = $wsAmsin.Cells.Item(1, 50)
$scratch.NumberFormat = $dateFmt
$scratch.Value = 1

$wsAmsin.Cells.Item(2, 2).NumberFormat = $dateFmt
$wsAmsin.Cells.Item(2, 2).Value = 44351.47052467593

Set-TextValue $wsAmsin.Cells.Item(3, 1) "2021-06-07"

$wsAmsin.Cells.Item(3, 2).NumberFormat = $dateFmt
$wsAmsin.Cells.Item(3, 2).Value = 44354.78134703704

$wsAmsin.Cells.Item(3, 3).Value = "regression_145final"
Set-TextValue $wsAmsin.Cells.Item(3, 4) "96"
$wsAmsin.Cells.Item(3, 5).Value = 96
$wsAmsin.Cells.Item(3, 6).Value = 0
$wsAmsin.Cells.Item(3, 7).Value = 2.809260433333333

$scratch.EntireColumn.Delete()

# ---------------------------------------------------------------------------
# BETA sheet: append the 2021-06-08 beta run as row 2.
# ---------------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Set-TextValue $wsBeta.Cells.Item(2, 1) "2021-06-08"

$wsBeta.Cells.Item(2, 2).NumberFormat = $dateFmt
$wsBeta.Cells.Item(2, 2).Value = 44355.69176928625

$wsBeta.Cells.Item(2, 3).Value = "145_beta"
Set-TextValue $wsBeta.Cells.Item(2, 4) "96"
$wsBeta.Cells.Item(2, 5).Value = 96
$wsBeta.Cells.Item(2, 6).Value = 0
$wsBeta.Cells.Item(2, 7).Value = 2.514197133333334
